# Auto-generated Word COM-interop edit script
$d = $word.ActiveDocument

# Smlouva 1 - Zamestnanec: phone -> tel.: +420 AMOUNT_1
$p1Find = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String("W1tQSE9ORV8xXV0="))
$p1Repl = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String("dGVsLjogKzQyMCBbW0FNT1VOVF8xXV0="))
$ok_p1 = $d.Content.Find.Execute($p1Find, $false, $false, $false, $false, $false, $true, 1, $false, $p1Repl, 2)
if (-not $ok_p1) { Write-Host "WARN: replace failed for p1" }

# Smlouva 1 - Pozice/Nastup/Mzda paragraph: merge bold-label runs into
# a single run without bold, collapsing 'Mzda: 55 000 Kc mesicne' into
# 'Mzda: [[AMOUNT_2]] mesicne'
$pozFindAnchor = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String("TXpkYTo="))
$pozPara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.Contains($pozFindAnchor) -and $para.Range.Text.Contains("Pozice:")) {
        $pozPara = $para
        break
    }
}
if ($pozPara -eq $null) { Write-Host "WARN: Pozice/Nastup/Mzda paragraph not found" }

# Clear bold on the paragraph's text content (exclude trailing paragraph mark)
$pozStart = $pozPara.Range.Start
$pozEnd = $pozPara.Range.End - 1
$pozContentRange = $d.Range($pozStart, $pozEnd)
$pozContentRange.Font.Bold = 0

$pozFindTxt = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String("UG96aWNlOiptxJtzw63EjW7Emw=="))
$pozReplTxt = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String("UG96aWNlOiBTb2Z0d2Fyb3bDvSBpbsW+ZW7DvXIg4oCUIHNlbmlvcgtOw6FzdHVwOiAxLiA2LiAyMDI1C016ZGE6IFtbQU1PVU5UXzJdXSBtxJtzw63EjW7Emw=="))
$pozRange = $pozPara.Range
$ok_poz = $pozRange.Find.Execute($pozFindTxt, $false, $false, $true, $false, $false, $true, 1, $false, $pozReplTxt, 2)
if (-not $ok_poz) { Write-Host "WARN: replace failed for Pozice/Nastup/Mzda" }

# Smlouva 2 - Pronajimatel: phone -> tel.: AMOUNT_3
$p3Find = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String("W1tQSE9ORV8yXV0="))
$p3Repl = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String("dGVsLjogW1tBTU9VTlRfM11d"))
$ok_p3 = $d.Content.Find.Execute($p3Find, $false, $false, $false, $false, $false, $true, 1, $false, $p3Repl, 2)
if (-not $ok_p3) { Write-Host "WARN: replace failed for p3" }

# Smlouva 2 - Najemce: phone -> tel.: AMOUNT_4
$p4Find = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String("W1tQSE9ORV8zXV0="))
$p4Repl = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String("dGVsLjogW1tBTU9VTlRfNF1d"))
$ok_p4 = $d.Content.Find.Execute($p4Find, $false, $false, $false, $false, $false, $true, 1, $false, $p4Repl, 2)
if (-not $ok_p4) { Write-Host "WARN: replace failed for p4" }

# Smlouva 2 - Najemne/Kauce amounts -> AMOUNT_5 / AMOUNT_6
$p5Find = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String("TsOhamVtbsOpOiAxNCA1MDAgS8SNIG3Em3PDrcSNbsSbLiBLYXVjZTogMjkgMDAwIEvEjS4="))
$p5Repl = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String("TsOhamVtbsOpOiBbW0FNT1VOVF81XV0gbcSbc8OtxI1uxJsuIEthdWNlOiBbW0FNT1VOVF82XV0u"))
$ok_p5 = $d.Content.Find.Execute($p5Find, $false, $false, $false, $false, $false, $true, 1, $false, $p5Repl, 2)
if (-not $ok_p5) { Write-Host "WARN: replace failed for p5" }

# Smlouva 3 - Prodavajici: phone -> tel.: AMOUNT_7
$p6Find = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String("W1tQSE9ORV80XV0="))
$p6Repl = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String("dGVsLjogW1tBTU9VTlRfN11d"))
$ok_p6 = $d.Content.Find.Execute($p6Find, $false, $false, $false, $false, $false, $true, 1, $false, $p6Repl, 2)
if (-not $ok_p6) { Write-Host "WARN: replace failed for p6" }

# Smlouva 3 - Kupujici: phone -> tel.: AMOUNT_8
$p7Find = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String("W1tQSE9ORV81XV0="))
$p7Repl = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String("dGVsLjogW1tBTU9VTlRfOF1d"))
$ok_p7 = $d.Content.Find.Execute($p7Find, $false, $false, $false, $false, $false, $true, 1, $false, $p7Repl, 2)
if (-not $ok_p7) { Write-Host "WARN: replace failed for p7" }

# Smlouva 3 - car price -> AMOUNT_9
$p8Find = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String("Y2VuYTogMTIwIDAwMCBLxI0="))
$p8Repl = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String("W1tBTU9VTlRfOV1dIEvEjQ=="))
$ok_p8 = $d.Content.Find.Execute($p8Find, $false, $false, $false, $false, $false, $true, 1, $false, $p8Repl, 2)
if (-not $ok_p8) { Write-Host "WARN: replace failed for p8" }

# Smlouva 4 - Objednatel: phone -> tel.: AMOUNT_10
$p9Find = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String("W1tQSE9ORV82XV0="))
$p9Repl = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String("dGVsLjogW1tBTU9VTlRfMTBdXQ=="))
$ok_p9 = $d.Content.Find.Execute($p9Find, $false, $false, $false, $false, $false, $true, 1, $false, $p9Repl, 2)
if (-not $ok_p9) { Write-Host "WARN: replace failed for p9" }

# Smlouva 4 - Predmet dila paragraph: merge the two runs spanning the
# <w:br/> into one, and turn 'Cena: 350 000 Kc' into '[[AMOUNT_11]] Kc'
$dilaFindTxt = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String("UMWZZWRtxJt0IGTDrWxhOipUZXJtw61uIGRva29uxI1lbsOtOiAzMC4gOS4gMjAyNS4="))
$dilaReplTxt = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String("UMWZZWRtxJt0IGTDrWxhOiBLb21wbGV4bsOtIHJla29uc3RydWtjZSBvcmRpbmFjZSDigJQgcG9kbGFoeSwgZWxla3RyaWthLCBtYWxiYS4LW1tBTU9VTlRfMTFdXSBLxI0gYmV6IERQSC4gVGVybcOtbiBkb2tvbsSNZW7DrTogMzAuIDkuIDIwMjUu"))
$ok_dila = $d.Content.Find.Execute($dilaFindTxt, $false, $false, $true, $false, $false, $true, 1, $false, $dilaReplTxt, 2)
if (-not $ok_dila) { Write-Host "WARN: replace failed for Predmet dila" }

# Smlouva 5 - Strana A: phone -> tel.: AMOUNT_12
$p11Find = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String("W1tQSE9ORV83XV0="))
$p11Repl = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String("dGVsLjogW1tBTU9VTlRfMTJdXQ=="))
$ok_p11 = $d.Content.Find.Execute($p11Find, $false, $false, $false, $false, $false, $true, 1, $false, $p11Repl, 2)
if (-not $ok_p11) { Write-Host "WARN: replace failed for p11" }

# Smlouva 5 - Kontaktni osoba: phone -> tel.: AMOUNT_13
$p12Find = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String("W1tQSE9ORV84XV0="))
$p12Repl = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String("dGVsLjogW1tBTU9VTlRfMTNdXQ=="))
$ok_p12 = $d.Content.Find.Execute($p12Find, $false, $false, $false, $false, $false, $true, 1, $false, $p12Repl, 2)
if (-not $ok_p12) { Write-Host "WARN: replace failed for p12" }

# Smlouva 6 - Pracovnik: phone -> tel.: +420 AMOUNT_14
$p13Find = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String("W1tQSE9ORV85XV0="))
$p13Repl = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String("dGVsLjogKzQyMCBbW0FNT1VOVF8xNF1d"))
$ok_p13 = $d.Content.Find.Execute($p13Find, $false, $false, $false, $false, $false, $true, 1, $false, $p13Repl, 2)
if (-not $ok_p13) { Write-Host "WARN: replace failed for p13" }

# Smlouva 7 - Veritel: phone -> tel.: AMOUNT_15
$p14Find = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String("W1tQSE9ORV8xMF1d"))
$p14Repl = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String("dGVsLjogW1tBTU9VTlRfMTVdXQ=="))
$ok_p14 = $d.Content.Find.Execute($p14Find, $false, $false, $false, $false, $false, $true, 1, $false, $p14Repl, 2)
if (-not $ok_p14) { Write-Host "WARN: replace failed for p14" }

# Smlouva 7 - pujcka amount -> AMOUNT_16
$p15Find = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String("UMWvasSNa2EgMTUwIDAwMCBLxI0="))
$p15Repl = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String("UMWvasSNa2EgW1tBTU9VTlRfMTZdXQ=="))
$ok_p15 = $d.Content.Find.Execute($p15Find, $false, $false, $false, $false, $false, $true, 1, $false, $p15Repl, 2)
if (-not $ok_p15) { Write-Host "WARN: replace failed for p15" }

# Smlouva 7 - Svedek: phone -> tel.: AMOUNT_17
$p16Find = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String("W1tQSE9ORV8xMV1d"))
$p16Repl = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String("dGVsLjogW1tBTU9VTlRfMTddXQ=="))
$ok_p16 = $d.Content.Find.Execute($p16Find, $false, $false, $false, $false, $false, $true, 1, $false, $p16Repl, 2)
if (-not $ok_p16) { Write-Host "WARN: replace failed for p16" }

# Smlouva 9 - kitchen price -> AMOUNT_18
$p17Find = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String("Y2VuYSAxOCAwMDAgS8SN"))
$p17Repl = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String("W1tBTU9VTlRfMThdXSBLxI0="))
$ok_p17 = $d.Content.Find.Execute($p17Find, $false, $false, $false, $false, $false, $true, 1, $false, $p17Repl, 2)
if (-not $ok_p17) { Write-Host "WARN: replace failed for p17" }

# Smlouva 10 - Uzivatel: -> USERNAME_1
$p18Find = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String("VcW+aXZhdGVsOg=="))
$p18Repl = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String("W1tVU0VSTkFNRV8xXV0="))
$ok_p18 = $d.Content.Find.Execute($p18Find, $false, $false, $false, $false, $false, $true, 1, $false, $p18Repl, 2)
if (-not $ok_p18) { Write-Host "WARN: replace failed for p18" }

Write-Host "DONE"
